# Added loading screen and new gamemodes
# Inserts new translation rows (login/nickname screen + quiz gameplay
# strings) into the Translations sheet, pushing the two existing rows
# (ANSWER100QUESTIONS / LOSE1TIME) further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 19 new blank rows before row 12 -------------------------------
# This shifts the existing rows 12-13 (ANSWER100QUESTIONS, LOSE1TIME) down
# to rows 31-32, making room for the new strings.
$ws.Range("A12:A30").EntireRow.Insert()

# --- Fill in the new rows (12-30) with Key / Polish / English values -----
$newRows = @(
    @("EMAIL",             "E-mail",                 "E-mail"),
    @("PASSWORD",          "Hasło",                  "Password"),
    @("NICKNAME",          "Nazwa",                   "Nickname"),
    @("REMEMBER",          "Zapamiętaj",              "Remember"),
    @("LOGIN",             "Login",                   "Login"),
    @("REGISTERANDLOGIN",  "Zarejestruj i zaloguj",   "Register and login"),
    @("RETURN",            "Powrót",                  "Return"),
    @("SEARCH",            "Szukaj…",                 "Search…"),
    @("SCORE",             "Wynik",                   "Score"),
    @("EASY",               "Łatwy",                  "Easy"),
    @("MEDIUM",            "Średni",                  "Medium"),
    @("HARD",              "Trudny",                  "Hard"),
    @("QUESTION",          "Pytanie",                  "Question"),
    @("TIME",              "Czas",                     "Time"),
    @("CORRECTANSWER",     "Poprawne odpowiedzi",      "Correct answers"),
    @("INCORRECTANSWER",   "Niepoprawne odpowiedzi",   "Incorrect answers"),
    @("TIMEISUP",          "Czas minał",                "Time is up"),
    @("PLAYAGAIN",         "Zagraj jeszcze raz",        "Play again"),
    @("REMAININGTIME",     "Pozostały czas",            "Remaining time")
)

$row = 12
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# --- Update the view: scroll position and active cell selection ----------
$ws.Range("A15").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
